$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from H1 into new header cells I1, J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I and J data columns
$iVals = @(
    6,
    4,
    5,
    2,
    6,
    6,
    6,
    8,
    8,
    7,
    10,
    7,
    7,
    7,
    7,
    8,
    8,
    7,
    7,
    7,
    7,
    7,
    6,
    6,
    9,
    7,
    7,
    10,
    8,
    9,
    8,
    8,
    8,
    9,
    9,
    7,
    9,
    7,
    9,
    9,
    8,
    8,
    7,
    8,
    7,
    6,
    6,
    8,
    7,
    7,
    8,
    7,
    10,
    5,
    9,
    6,
    6,
    6,
    3,
    7,
    7,
    7,
    1,
    7
)
$jVals = @(
    6,
    5,
    6,
    3,
    6,
    7,
    6,
    8,
    8,
    7,
    10,
    7,
    8,
    8,
    8,
    8,
    8,
    7,
    7,
    7,
    7,
    8,
    6,
    6,
    9,
    7,
    7,
    10,
    8,
    9,
    9,
    8,
    8,
    9,
    9,
    9,
    9,
    7,
    9,
    9,
    8,
    8,
    7,
    9,
    7,
    7,
    7,
    8,
    7,
    7,
    9,
    7,
    11,
    7,
    9,
    6,
    7,
    7,
    5,
    9,
    8,
    7,
    1,
    7
)

$startRow = 2
for ($k = 0; $k -lt $iVals.Length; $k++) {
    $r = $startRow + $k
    $ws.Cells.Item($r, 9).Value = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}
